# Auto-generated: refresh market-price derived columns (H-N) across all Leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "M41" = -5290776
    "K41" = 5291216
    "H41" = 3831647.2
    "L41" = 280.25
    "I41" = 5291216
    "J41" = 280.25
    "N41" = -1160.25
    "H132" = 336029.22
    "L132" = 220577.25
    "I132" = 406030.12
    "J132" = 73525.75
    "N132" = -225637.25
    "M132" = -1215560.36
    "K132" = 1218090.36
    "N134" = -60140
    "H134" = 50000
    "J134" = 50000
    "L134" = 50000
    "H135" = 1780.381
    "L135" = 20880
    "I135" = 1611.75
    "J135" = 2320
    "N135" = -25950
    "M135" = -11970.75
    "K135" = 14505.75
    "L136" = 56666.668
    "J136" = 56666.668
    "N136" = -66866.66800000001
    "H136" = 56666.668
    "H139" = 49800
    "L139" = 49800
    "J139" = 49800
    "N139" = -60080
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H32" = 13036.784
    "I32" = 1490.6385
    "M32" = -1203.6385
    "K32" = 1490.6385
    "H45" = 960
    "L45" = 1000
    "I45" = 950
    "J45" = 1000
    "N45" = -1754
    "M45" = -573
    "K45" = 950
    "H122" = 2137.4
    "I122" = 1802.2
    "J122" = 2472.6
    "L122" = 7417.799999999999
    "N122" = -12317.8
    "M122" = -2956.6
    "K122" = 5406.6
    "L133" = 49505.668
    "J133" = 49505.668
    "N133" = -54565.668
    "H133" = 49505.668
    "N134" = -70140
    "H134" = 60000
    "J134" = 60000
    "L134" = 60000
    "H135" = 31608.285
    "L135" = 31608.285
    "J135" = 31608.285
    "N135" = -41748.285
    "H137" = 44999.5
    "L137" = 44999.5
    "J137" = 44999.5
    "N137" = -55199.5
    "H139" = 43665.8
    "L139" = 43665.8
    "J139" = 43665.8
    "N139" = -53945.8
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "L59" = 47280
    "J59" = 47280
    "N59" = -48974
    "H59" = 42733.332
    "H132" = 21779
    "L132" = 21779
    "J132" = 21779
    "N132" = -31899
    "L133" = 39334.75
    "J133" = 39334.75
    "N133" = -49454.75
    "H133" = 39334.75
    "H135" = 57694.445
    "L135" = 57694.445
    "J135" = 57694.445
    "N135" = -67834.44500000001
    "H137" = 46183.332
    "L137" = 46183.332
    "J137" = 46183.332
    "N137" = -56383.332
    "H138" = 0
    "L138" = 0
    "J138" = 0
    "N140" = -57620
    "H140" = 47260
    "L140" = 47260
    "J140" = 47260
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "M2" = -891
    "K2" = 1004
    "H2" = 47004.668
    "L2" = 70005
    "I2" = 1004
    "J2" = 70005
    "N2" = -70231
    "H94" = 1450.2858
    "I94" = 902.8889
    "M94" = -451.8889
    "K94" = 902.8889
    "H137" = 42468
    "L137" = 46075.555
    "J137" = 46075.555
    "N137" = -56275.555
    "N140" = -57942.5
    "H140" = 47582.5
    "L140" = 47582.5
    "J140" = 47582.5
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H4" = 31126.545
    "L4" = 146573.145
    "I4" = 97
    "J4" = 48857.715
    "N4" = -146797.145
    "M4" = -179
    "K4" = 291
    "H51" = 875
    "J51" = 800
    "L51" = 2400
    "N51" = -3320
    "H55" = 24834.334
    "L55" = 83064.375
    "I55" = 2004
    "J55" = 27688.125
    "N55" = -83418.375
    "M55" = -5835
    "K55" = 6012
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H5" = 433.33334
    "I46" = 12080
    "L46" = 23200
    "J46" = 23200
    "N46" = -23512
    "M46" = -11924
    "K46" = 12080
    "H46" = 20420
    "H113" = 1700
    "L113" = 1933.3334
    "I113" = 1583.3334
    "J113" = 1933.3334
    "N113" = -6273.3334
    "M113" = 586.6666
    "K113" = 1583.3334
    "H137" = 55000
    "L137" = 55000
    "J137" = 55000
    "N137" = -65200
    "H138" = 59574.875
    "L138" = 59574.875
    "J138" = 59574.875
    "N138" = -69854.875
    "H139" = 36163
    "L139" = 36163
    "J139" = 36163
    "N139" = -46443
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H55" = 491.7143
    "L55" = 528.4
    "I55" = 400
    "J55" = 528.4
    "N55" = -874.4
    "M55" = -227
    "K55" = 400
    "H61" = 4347.913
    "L61" = 2999
    "J61" = 2999
    "N61" = -3403
    "H113" = 4347.913
    "L113" = 2999
    "J113" = 2999
    "N113" = -7339
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H81" = 1668043.4
    "L81" = 5500
    "I81" = 2001102
    "J81" = 2750
    "N81" = -7622
    "M81" = -4001143
    "K81" = 4002204
    "H84" = 1668043.4
    "I84" = 2001102
    "J84" = 2750
    "L84" = 27500
    "N84" = -38108
    "M84" = -20005716
    "K84" = 20011020
    "H132" = 12198463
    "L132" = 7250.625
    "I132" = 20003932
    "J132" = 2416.875
    "N132" = -12310.625
    "M132" = -60009266
    "K132" = 60011796
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
